$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure affected cells keep a Text number format so values
# such as "0.9998" or "1.000" are not re-interpreted as numbers/dates.
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '29.215.19'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  -0.69%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.861.79'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  -0.80%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.9998'
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '0.7149'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  -0.33%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '240.46'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  +0.16%  '
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  +0.05%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3088'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  -0.30%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07702'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  -1.73%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '24.98'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  +0.75%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.08311'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  +0.66%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.895.57'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  +1.35%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.7175'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  -1.25%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.216'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  -1.26%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '90.89'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  -0.37%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '29.254.86'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  -0.53%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '5.993'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  +1.41%  '
$ws.Range('B18').NumberFormat = '@'
$ws.Range('B18').Value = 'BitcoinCash'
$ws.Range('C18').NumberFormat = '@'
$ws.Range('C18').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '243.52'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  -0.64%  '
$ws.Range('B19').NumberFormat = '@'
$ws.Range('B19').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C19').NumberFormat = '@'
$ws.Range('C19').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '2.159.69'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  +1.91%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.000007808'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  -1.09%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '13.16'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  -0.97%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.9998'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  -0.02%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '7.965'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  +0.80%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '0.9999'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  +0.00%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.1615'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  +3.25%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '162.84'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  -0.63%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '8.904'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  -1.23%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '18.59'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  +1.35%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.350'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  -0.50%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '4.443'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  +1.11%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.497'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  +0.87%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.248'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  +2.52%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.05181'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  -1.89%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.7973'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  +10.37%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.932'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  -0.16%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.172'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  -2.35%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.682'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  +0.24%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.01858'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  -0.54%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.691'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  -0.98%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.174.32'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  -4.70%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '6.213'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  +2.26%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.9026'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  -0.80%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '72.96'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  +0.02%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.9992'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  -0.06%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.054.60'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  +1.91%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '102.17'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  -1.50%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.5197'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  -2.66%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.782'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  +1.32%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '9.375'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  +1.30%  '
$ws.Range('B50').NumberFormat = '@'
$ws.Range('B50').Value = 'BabyDogeCoin'
$ws.Range('C50').NumberFormat = '@'
$ws.Range('C50').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.00000000120'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  -0.98%  '
$ws.Range('B51').NumberFormat = '@'
$ws.Range('B51').Value = 'Aptos'
$ws.Range('C51').NumberFormat = '@'
$ws.Range('C51').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '7.078'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  +0.23%  '
